# Add files via upload
# Populate the leaderboard data on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("Jack",     4, 4, 1),
    @("Lance ",   3, 3, 1),
    @("Connor",   2, 1, 2),
    @("Phillip",  1, 2, 0.5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
}

$null = $ws.Range("D5").Select()
